# Apply "photoprotection catalogue" update to the Sunscreens sheet.
#  - Row 2 (Tesco £8.50 product): price cell becomes a text label, UVA label
#    swaps to "UVB", several lab-measured figures change.
#  - Row 3 (Boots £14 product): price cell becomes a text label, UVA label
#    swaps to "UVB UVA", volume + lab figures change.
#  - Row 4 (Specialist £22 product): price cell becomes a text label,
#    volume + lab figures change.
#  - Rows 5 and 6 are brand-new "Mystery" entries (Sunscreen A / Sunscreen B).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sunscreens")

# ---- Row 2 ---------------------------------------------------------------
$ws.Range("C2").Value = "£8.50 sunscreen"
$ws.Range("D2").Value = "UVB"
$ws.Range("T2").Value = 80
$ws.Range("U2").Value = 0.01
$ws.Range("V2").Value = 0.005
$ws.Range("W2").Value = 5

# ---- Row 3 ---------------------------------------------------------------
$ws.Range("C3").Value = "£14 sunscreen"
$ws.Range("D3").Value = "UVB UVA"
$ws.Range("K3").Value = 100
$ws.Range("T3").Value = 100
$ws.Range("U3").Value = 0.02
$ws.Range("W3").Value = 40

# ---- Row 4 ---------------------------------------------------------------
$ws.Range("C4").Value = "£22 sunscreen"
$ws.Range("K4").Value = 50
$ws.Range("U4").Value = 0.5
$ws.Range("V4").Value = 0.4
$ws.Range("W4").Value = 25

# ---- Row 5 (new) -----------------------------------------------------------
$ws.Range("C5").Value = "Mystery "
$ws.Range("D5").Value = "Sunscreen A"
$ws.Range("K5").Value = 200
$ws.Range("L5").Value = 8.5
$ws.Range("M5").Formula = "=L5/K5"
$ws.Range("T5").Value = 90
$ws.Range("U5").Value = 0.01
$ws.Range("V5").Value = 0.005
$ws.Range("W5").Value = 40

# ---- Row 6 (new) -----------------------------------------------------------
$ws.Range("C6").Value = "Mystery "
$ws.Range("D6").Value = "Sunscreen B"
$ws.Range("K6").Value = 50
$ws.Range("L6").Value = 25
$ws.Range("M6").Formula = "=L6/K6"
$ws.Range("T6").Value = 40
$ws.Range("U6").Value = 0.005
$ws.Range("V6").Value = 0.005
$ws.Range("W6").Value = 1

# ---- View state: scroll so column D is left-most and X6 is selected -------
$ws.Activate()
$ws.Range("X6").Select()
